$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1047491833"
$ws.Range("D16").Value = "JOSE HERNAN ROMERO ATENCIO"
$ws.Range("E16").Value = "2502"
$ws.Range("F16").Value = 32933
$ws.Range("G16").Value = 1300000

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1047491833"
$ws.Range("D17").Value = "JOSE HERNAN ROMERO ATENCIO"
$ws.Range("E17").Value = "2501"
$ws.Range("F17").Value = 52000
$ws.Range("G17").Value = 1300000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1047491833"
$ws.Range("D18").Value = "JOSE HERNAN ROMERO ATENCIO"
$ws.Range("E18").Value = "2412"
$ws.Range("F18").Value = 52000
$ws.Range("G18").Value = 1300000

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1047491833"
$ws.Range("D19").Value = "JOSE HERNAN ROMERO ATENCIO"
$ws.Range("E19").Value = "2411"
$ws.Range("F19").Value = 52000
$ws.Range("G19").Value = 1300000

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1047491833"
$ws.Range("D20").Value = "JOSE HERNAN ROMERO ATENCIO"
$ws.Range("E20").Value = "2410"
$ws.Range("F20").Value = 52000
$ws.Range("G20").Value = 1300000

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1047491833"
$ws.Range("D21").Value = "JOSE HERNAN ROMERO ATENCIO"
$ws.Range("E21").Value = "2409"
$ws.Range("F21").Value = 52000
$ws.Range("G21").Value = 1300000

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1047491833"
$ws.Range("D22").Value = "JOSE HERNAN ROMERO ATENCIO"
$ws.Range("E22").Value = "2408"
$ws.Range("F22").Value = 52000
$ws.Range("G22").Value = 1300000

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1047491833"
$ws.Range("D23").Value = "JOSE HERNAN ROMERO ATENCIO"
$ws.Range("E23").Value = "2407"
$ws.Range("F23").Value = 52000
$ws.Range("G23").Value = 1300000

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1047491833"
$ws.Range("D24").Value = "JOSE HERNAN ROMERO ATENCIO"
$ws.Range("E24").Value = "2406"
$ws.Range("F24").Value = 52000
$ws.Range("G24").Value = 1300000

$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "1047491833"
$ws.Range("D25").Value = "JOSE HERNAN ROMERO ATENCIO"
$ws.Range("E25").Value = "2405"
$ws.Range("F25").Value = 52000
$ws.Range("G25").Value = 1300000

$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "1047491833"
$ws.Range("D26").Value = "JOSE HERNAN ROMERO ATENCIO"
$ws.Range("E26").Value = "2404"
$ws.Range("F26").Value = 52000
$ws.Range("G26").Value = 1300000

$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "1047491833"
$ws.Range("D27").Value = "JOSE HERNAN ROMERO ATENCIO"
$ws.Range("E27").Value = "2403"
$ws.Range("F27").Value = 45066
$ws.Range("G27").Value = 1300000

$ws.Range("B28").Value = "CC"
$ws.Range("C28").Value = "7919980"
$ws.Range("D28").Value = "ENRIQUE CARLOS CASTRO SALCEDO"
$ws.Range("E28").Value = "2502"
$ws.Range("F28").Value = 32933
$ws.Range("G28").Value = 1300000

$ws.Range("B29").Value = "CC"
$ws.Range("C29").Value = "7919980"
$ws.Range("D29").Value = "ENRIQUE CARLOS CASTRO SALCEDO"
$ws.Range("E29").Value = "2501"
$ws.Range("F29").Value = 52000
$ws.Range("G29").Value = 1300000

$ws.Range("B30").Value = "CC"
$ws.Range("C30").Value = "7919980"
$ws.Range("D30").Value = "ENRIQUE CARLOS CASTRO SALCEDO"
$ws.Range("E30").Value = "2412"
$ws.Range("F30").Value = 52000
$ws.Range("G30").Value = 1300000

$ws.Range("B31").Value = "CC"
$ws.Range("C31").Value = "7919980"
$ws.Range("D31").Value = "ENRIQUE CARLOS CASTRO SALCEDO"
$ws.Range("E31").Value = "2411"
$ws.Range("F31").Value = 52000
$ws.Range("G31").Value = 1300000

$ws.Range("B32").Value = "CC"
$ws.Range("C32").Value = "7919980"
$ws.Range("D32").Value = "ENRIQUE CARLOS CASTRO SALCEDO"
$ws.Range("E32").Value = "2410"
$ws.Range("F32").Value = 52000
$ws.Range("G32").Value = 1300000

$ws.Range("B33").Value = "CC"
$ws.Range("C33").Value = "7919980"
$ws.Range("D33").Value = "ENRIQUE CARLOS CASTRO SALCEDO"
$ws.Range("E33").Value = "2409"
$ws.Range("F33").Value = 52000
$ws.Range("G33").Value = 1300000

$ws.Range("B34").Value = "CC"
$ws.Range("C34").Value = "7919980"
$ws.Range("D34").Value = "ENRIQUE CARLOS CASTRO SALCEDO"
$ws.Range("E34").Value = "2408"
$ws.Range("F34").Value = 52000
$ws.Range("G34").Value = 1300000

$ws.Range("B35").Value = "CC"
$ws.Range("C35").Value = "1050968241"
$ws.Range("D35").Value = "KEYLA PAOLA PATERNINA TORCUATO"
$ws.Range("E35").Value = "2409"
$ws.Range("F35").Value = 52000
$ws.Range("G35").Value = 1300000

$ws.Range("B36").Value = "CC"
$ws.Range("C36").Value = "1050968241"
$ws.Range("D36").Value = "KEYLA PAOLA PATERNINA TORCUATO"
$ws.Range("E36").Value = "2408"
$ws.Range("F36").Value = 52000
$ws.Range("G36").Value = 1300000

$ws.Range("B37").Value = "CC"
$ws.Range("C37").Value = "1046274326"
$ws.Range("D37").Value = "CARLOS MARIO HERNANDEZ MORALES"
$ws.Range("E37").Value = "2502"
$ws.Range("F37").Value = 32933
$ws.Range("G37").Value = 828116

$ws.Range("B38").Value = "CC"
$ws.Range("C38").Value = "1046274326"
$ws.Range("D38").Value = "CARLOS MARIO HERNANDEZ MORALES"
$ws.Range("E38").Value = "2501"
$ws.Range("F38").Value = 52000
$ws.Range("G38").Value = 828116

$ws.Range("B39").Value = "CC"
$ws.Range("C39").Value = "1046274326"
$ws.Range("D39").Value = "CARLOS MARIO HERNANDEZ MORALES"
$ws.Range("E39").Value = "2412"
$ws.Range("F39").Value = 52000
$ws.Range("G39").Value = 828116

$ws.Range("B40").Value = "CC"
$ws.Range("C40").Value = "1046274326"
$ws.Range("D40").Value = "CARLOS MARIO HERNANDEZ MORALES"
$ws.Range("E40").Value = "2411"
$ws.Range("F40").Value = 52000
$ws.Range("G40").Value = 828116

$ws.Range("B41").Value = "CC"
$ws.Range("C41").Value = "1046274326"
$ws.Range("D41").Value = "CARLOS MARIO HERNANDEZ MORALES"
$ws.Range("E41").Value = "2410"
$ws.Range("F41").Value = 52000
$ws.Range("G41").Value = 828116

$ws.Range("B42").Value = "CC"
$ws.Range("C42").Value = "1046274326"
$ws.Range("D42").Value = "CARLOS MARIO HERNANDEZ MORALES"
$ws.Range("E42").Value = "2409"
$ws.Range("F42").Value = 52000
$ws.Range("G42").Value = 828116

$ws.Range("B43").Value = "CC"
$ws.Range("C43").Value = "1046274326"
$ws.Range("D43").Value = "CARLOS MARIO HERNANDEZ MORALES"
$ws.Range("E43").Value = "2408"
$ws.Range("F43").Value = 52000
$ws.Range("G43").Value = 828116

$ws.Range("B44").Value = "CC"
$ws.Range("C44").Value = "1047510512"
$ws.Range("D44").Value = "GERMAN DAVID MAZO SIERRA"
$ws.Range("E44").Value = "2409"
$ws.Range("F44").Value = 52000
$ws.Range("G44").Value = 1300000

$ws.Range("B45").Value = "CC"
$ws.Range("C45").Value = "1047510512"
$ws.Range("D45").Value = "GERMAN DAVID MAZO SIERRA"
$ws.Range("E45").Value = "2408"
$ws.Range("F45").Value = 52000
$ws.Range("G45").Value = 1300000

$ws.Range("B46").Value = "CC"
$ws.Range("C46").Value = "1003714323"
$ws.Range("D46").Value = "MARIA CLAUDIA MADRID ROSSO"
$ws.Range("E46").Value = "2409"
$ws.Range("F46").Value = 52000
$ws.Range("G46").Value = 1300000

$ws.Range("B47").Value = "PPT"
$ws.Range("C47").Value = "5943946"
$ws.Range("D47").Value = "GLEIDIMAR DEL CARMEN BETANCOURTH CAICEDO"
$ws.Range("E47").Value = "2409"
$ws.Range("F47").Value = 52000
$ws.Range("G47").Value = 1300000

$ws.Range("B48").Value = "PPT"
$ws.Range("C48").Value = "5943946"
$ws.Range("D48").Value = "GLEIDIMAR DEL CARMEN BETANCOURTH CAICEDO"
$ws.Range("E48").Value = "2408"
$ws.Range("F48").Value = 52000
$ws.Range("G48").Value = 1300000

